$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.706.92"
$ws.Range("E2").Value = "  -0.11%  "

# Row 3
$ws.Range("D3").Value = "2.529.13"
$ws.Range("E3").Value = "  +0.19%  "

# Row 4
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").Value = "315.57"
$ws.Range("E5").Value = "  -0.81%  "

# Row 6
$ws.Range("D6").Value = "95.52"
$ws.Range("E6").Value = "  -0.58%  "

# Row 7
$ws.Range("D7").Value = "0.572"
$ws.Range("E7").Value = "  -1.84%  "

# Row 8
$ws.Range("E8").Value = "  -0.02%  "

# Row 9
$ws.Range("D9").Value = "0.531"
$ws.Range("E9").Value = "  -1.53%  "

# Row 10
$ws.Range("D10").Value = "35.76"
$ws.Range("E10").Value = "  -1.46%  "

# Row 11
$ws.Range("D11").Value = "0.0804"
$ws.Range("E11").Value = "  -1.03%  "

# Row 12
$ws.Range("D12").Value = "7.52"
$ws.Range("E12").Value = "  -1.07%  "

# Row 13
$ws.Range("E13").Value = "  -2.79%  "

# Row 14
$ws.Range("D14").Value = "2.919.13"
$ws.Range("E14").Value = "  +0.25%  "

# Row 15
$ws.Range("D15").Value = "2.574.11"
$ws.Range("E15").Value = "  +2.85%  "

# Row 16
$ws.Range("D16").Value = "15.11"
$ws.Range("E16").Value = "  -2.88%  "

# Row 17
$ws.Range("D17").Value = "0.846"
$ws.Range("E17").Value = "  -1.82%  "

# Row 18
$ws.Range("D18").Value = "42.812.17"
$ws.Range("E18").Value = "  +0.20%  "

# Row 19
$ws.Range("D19").Value = "6.85"
$ws.Range("E19").Value = "  +3.36%  "

# Row 20
$ws.Range("D20").Value = "12.79"
$ws.Range("E20").Value = "  -1.05%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0960"
$ws.Range("E21").Value = "  -1.31%  "

# Row 22
$ws.Range("D22").Value = "'69.60"
$ws.Range("E22").Value = "  -2.69%  "

# Row 23
$ws.Range("D23").Value = "251.18"
$ws.Range("E23").Value = "  -0.95%  "

# Row 24
$ws.Range("D24").Value = "2.94"
$ws.Range("E24").Value = "  -1.86%  "

# Row 25
$ws.Range("D25").Value = "2.06"
$ws.Range("E25").Value = "  +0.81%  "

# Row 26
$ws.Range("D26").Value = "26.45"
$ws.Range("E26").Value = "  -2.34%  "

# Row 27
$ws.Range("E27").Value = "  +0.01%  "

# Row 28
$ws.Range("D28").Value = "2.41"
$ws.Range("E28").Value = "  +1.06%  "

# Row 29
$ws.Range("D29").Value = "40.32"
$ws.Range("E29").Value = "  +5.12%  "

# Row 30
$ws.Range("D30").Value = "10.38"
$ws.Range("E30").Value = "  +2.30%  "

# Row 31
$ws.Range("D31").Value = "5.89"
$ws.Range("E31").Value = "  -0.64%  "

# Row 32
$ws.Range("D32").Value = "156.44"
$ws.Range("E32").Value = "  +0.92%  "

# Row 33
$ws.Range("D33").Value = "2.15"
$ws.Range("E33").Value = "  +2.58%  "

# Row 34
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "3.32"
$ws.Range("E34").Value = "  -0.06%  "

# Row 35
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "2.69"
$ws.Range("E35").Value = "  +2.75%  "

# Row 36
$ws.Range("D36").Value = "18.81"
$ws.Range("E36").Value = "  -5.69%  "

# Row 37
$ws.Range("D37").Value = "0.0775"
$ws.Range("E37").Value = "  -1.92%  "

# Row 38
$ws.Range("D38").Value = "0.111"
$ws.Range("E38").Value = "  -1.73%  "

# Row 39
$ws.Range("E39").Value = "  -1.69%  "

# Row 40
$ws.Range("D40").Value = "2.32"
$ws.Range("E40").Value = "  +10.96%  "

# Row 41
$ws.Range("D41").Value = "22.26"
$ws.Range("E41").Value = "  -8.71%  "

# Row 42
$ws.Range("D42").Value = "3.81"
$ws.Range("E42").Value = "  -1.28%  "

# Row 43
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  +0.27%  "

# Row 44
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "0.0302"
$ws.Range("E44").Value = "  -0.12%  "

# Row 45
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.029.28"
$ws.Range("E45").Value = "  -0.10%  "

# Row 46
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "3.24"
$ws.Range("E46").Value = "  -4.73%  "

# Row 47
$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D47").Value = "84.46"
$ws.Range("E47").Value = "  -0.41%  "

# Row 48
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "9.03"
$ws.Range("E48").Value = "  +0.76%  "

# Row 49
$ws.Range("D49").Value = "105.81"
$ws.Range("E49").Value = "  +3.78%  "

# Row 50
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").Value = "74.69"
$ws.Range("E50").Value = "  +0.44%  "

# Row 51
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.773.71"
$ws.Range("E51").Value = "  +0.26%  "

